$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels with proper capitalization
$ws.Range("A1").Value = "Gene"
$ws.Range("B1").Value = "Primer"
$ws.Range("C1").Value = "HKG"

# Add new column for technical replicates workflow
$ws.Range("C2").Value = "TRUE OR FALSE"

# Leave selection on the newly edited cell
$ws.Range("C2").Select()
